$d = $word.ActiveDocument

$replacements = @(
    @("94÷5=18, 4", "40÷5=8, 0"),
    @("60÷7=8, 4", "43÷7=6, 1"),
    @("58÷2=29, 0", "83÷8=10, 3"),
    @("37÷7=5, 2", "20÷4=5, 0"),
    @("16÷7=2, 2", "88÷8=11, 0"),
    @("70÷4=17, 2", "33÷6=5, 3"),
    @("56÷5=11, 1", "83÷6=13, 5"),
    @("41÷4=10, 1", "11÷6=1, 5"),
    @("13÷7=1, 6", "61÷6=10, 1"),
    @("13÷6=2, 1", "56÷8=7, 0"),
    @("17÷6=2, 5", "83÷9=9, 2"),
    @("83÷7=11, 6", "34÷5=6, 4"),
    @("50÷7=7, 1", "48÷2=24, 0"),
    @("58÷8=7, 2", "29÷9=3, 2"),
    @("46÷8=5, 6", "15÷5=3, 0"),
    @("25÷4=6, 1", "91÷3=30, 1"),
    @("76÷2=38, 0", "29÷3=9, 2"),
    @("45÷9=5, 0", "94÷9=10, 4"),
    @("60÷3=20, 0", "56÷6=9, 2"),
    @("24÷8=3, 0", "48÷9=5, 3"),
    @("83÷5=16, 3", "95÷4=23, 3"),
    @("25÷3=8, 1", "57÷8=7, 1"),
    @("39÷2=19, 1", "73÷6=12, 1"),
    @("63÷2=31, 1", "92÷2=46, 0"),
    @("31÷3=10, 1", "25÷5=5, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
